# Updates the cryptos price/volume table to the latest scraped snapshot.
# Numeric-looking "Price" strings (e.g. "258.58") are written with a leading
# apostrophe so Excel stores them as text (same as a user typing '258.58 into
# the cell) instead of silently coercing them to a Number; this preserves the
# exact text representation (trailing zeros, etc.) the diff expects.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.937.30'
$ws.Range('E2').Value = '  +3.54%  '
$ws.Range('D3').Value = '2.235.21'
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''258.58'
$ws.Range('E5').Value = '  +1.76%  '
$ws.Range('D6').Value = '''80.36'
$ws.Range('E6').Value = '  +8.95%  '
$ws.Range('D7').Value = '''0.625'
$ws.Range('E7').Value = '  +2.30%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '''0.600'
$ws.Range('E9').Value = '  +1.95%  '
$ws.Range('D10').Value = '''43.12'
$ws.Range('E10').Value = '  +7.39%  '
$ws.Range('D11').Value = '''0.0925'
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').Value = '''7.04'
$ws.Range('E12').Value = '  +3.30%  '
$ws.Range('E13').Value = '  +2.19%  '
$ws.Range('D14').Value = '2.576.18'
$ws.Range('E14').Value = '  +2.08%  '
$ws.Range('D15').Value = '''14.68'
$ws.Range('E15').Value = '  +2.50%  '
$ws.Range('D16').Value = '2.242.83'
$ws.Range('E16').Value = '  +2.62%  '
$ws.Range('D17').Value = '''0.789'
$ws.Range('E17').Value = '  +1.81%  '
$ws.Range('D18').Value = '43.864.22'
$ws.Range('E18').Value = '  +3.64%  '
$ws.Range('E19').Value = '  +1.55%  '
$ws.Range('D20').Value = '''71.19'
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').Value = '''6.04'
$ws.Range('E21').Value = '  +2.25%  '
$ws.Range('D22').Value = '''2.38'
$ws.Range('E22').Value = '  +8.55%  '
$ws.Range('D23').Value = '''232.58'
$ws.Range('E23').Value = '  +2.06%  '
$ws.Range('D24').Value = '''9.48'
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('D26').Value = '''10.82'
$ws.Range('E26').Value = '  +2.02%  '
$ws.Range('D27').Value = '''40.62'
$ws.Range('E27').Value = '  +10.08%  '
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('E30').Value = '  -0.50%  '
$ws.Range('D31').Value = '''172.26'
$ws.Range('E31').Value = '  +1.97%  '
$ws.Range('D32').Value = '''0.0890'
$ws.Range('E32').Value = '  +10.28%  '
$ws.Range('D33').Value = '''20.58'
$ws.Range('E33').Value = '  +2.63%  '
$ws.Range('D34').Value = '''5.31'
$ws.Range('E34').Value = '  +2.94%  '
$ws.Range('D35').Value = '''0.116'
$ws.Range('E35').Value = '  +8.20%  '
$ws.Range('E36').Value = '  +1.91%  '
$ws.Range('D37').Value = '''0.0368'
$ws.Range('E37').Value = '  +11.18%  '
$ws.Range('D38').Value = '''4.50'
$ws.Range('E38').Value = '  +4.14%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = '''3.03'
$ws.Range('E39').Value = '  +25.85%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '''13.02'
$ws.Range('E40').Value = '  +7.66%  '
$ws.Range('E41').Value = '  +3.02%  '
$ws.Range('D42').Value = '''63.11'
$ws.Range('E42').Value = '  +7.00%  '
$ws.Range('D43').Value = '''5.51'
$ws.Range('E43').Value = '  +6.17%  '
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('D45').Value = '''104.11'
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('D46').Value = '''8.49'
$ws.Range('E46').Value = '  +2.00%  '
$ws.Range('E47').Value = '  +1.37%  '
$ws.Range('D48').Value = '''1.12'
$ws.Range('E48').Value = '  +2.14%  '
$ws.Range('D49').Value = '''1.56'
$ws.Range('E49').Value = '  +28.39%  '
$ws.Range('E50').Value = '  -5.19%  '
$ws.Range('D51').Value = '''1.15'
$ws.Range('E51').Value = '  +2.23%  '
